$wb = $excel.ActiveWorkbook

# --- Sheet1 (存款): insert a new detail row after the header/total row ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2, 1).Value = 42
$ws1.Cells.Item(2, 2).Value = "中國信託商業銀行斗六分行"
$ws1.Cells.Item(2, 3).Value = "活期儲蓄存款"
$ws1.Cells.Item(2, 4).Value = "新臺幣"
$ws1.Cells.Item(2, 5).Value = "劉建國"
$ws1.Cells.Item(2, 6).Value = 46845

$ws1.Cells.Item(2, 1).Style = $ws1.Cells.Item(3, 1).Style
$ws1.Cells.Item(2, 2).Style = $ws1.Cells.Item(3, 2).Style
$ws1.Cells.Item(2, 3).Style = $ws1.Cells.Item(3, 3).Style
$ws1.Cells.Item(2, 4).Style = $ws1.Cells.Item(3, 4).Style
$ws1.Cells.Item(2, 5).Style = $ws1.Cells.Item(3, 5).Style
$ws1.Cells.Item(2, 6).Style = $ws1.Cells.Item(3, 6).Style

# --- Sheet2 (保險): add new detail row below the existing row ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = 81
$ws2.Cells.Item(2, 2).Value = "三商美邦人壽"
$ws2.Cells.Item(2, 3).Value = "雙喜臨門B型儲蓄險"
$ws2.Cells.Item(2, 4).Value = "劉建國"
$ws2.Cells.Item(2, 5).Value = "保險年齡至保險人四If四歲繳費期間六年年繳二十萬"

$ws2.Cells.Item(2, 1).Style = $ws1.Cells.Item(3, 1).Style
$ws2.Cells.Item(2, 2).Style = $ws1.Cells.Item(3, 2).Style
$ws2.Cells.Item(2, 3).Style = $ws1.Cells.Item(3, 3).Style
$ws2.Cells.Item(2, 4).Style = $ws1.Cells.Item(3, 4).Style
$ws2.Cells.Item(2, 5).Style = $ws1.Cells.Item(3, 5).Style

# --- Sheet3 (債務): add new detail row below the existing row ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 1).Value = 91
$ws3.Cells.Item(2, 2).Value = "貸款"
$ws3.Cells.Item(2, 3).Value = "劉建國"
$ws3.Cells.Item(2, 4).Value = "斗六市農會雲林縣斗六市民生路"
$ws3.Cells.Item(2, 5).Value = 5200000
$ws3.Cells.Item(2, 6).Value = "99年12月08日"
$ws3.Cells.Item(2, 7).Value = "代償債務"

$ws3.Cells.Item(2, 1).Style = $ws1.Cells.Item(3, 1).Style
$ws3.Cells.Item(2, 2).Style = $ws1.Cells.Item(3, 2).Style
$ws3.Cells.Item(2, 3).Style = $ws1.Cells.Item(3, 3).Style
$ws3.Cells.Item(2, 4).Style = $ws1.Cells.Item(3, 4).Style
$ws3.Cells.Item(2, 5).Style = $ws1.Cells.Item(3, 5).Style
$ws3.Cells.Item(2, 6).Style = $ws1.Cells.Item(3, 6).Style
$ws3.Cells.Item(2, 7).Style = $ws1.Cells.Item(3, 6).Style

$wb.Save()
